$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program dependencies")
Write-Host $ws.Range("A2").Formula
Write-Host $ws.Range("A1").Formula
$ws.Range("Z1").Value = 42
Write-Host $ws.Range("Z1").Formula
